# Updated the Web Inspect
#
# Adds two new locator sheets ("ImageLocators" and "HeadingLocators") after
# the existing "TableLocators" sheet, gives each the same Loc1..Loc7 header
# row used by every other *Locators sheet in this workbook, refreshes the
# selection on "TableLocators" (it is no longer the active tab), and leaves
# "HeadingLocators" as the new active sheet/tab with cell B2 selected.

$wb = $excel.ActiveWorkbook

$headers = @("Loc1", "Loc2", "Loc3", "Loc4", "Loc5", "Loc6", "Loc7")

# --- TableLocators: no longer the active tab; selection resets to the header row ---
$tableSheet = $wb.Worksheets.Item("TableLocators")
$tableSheet.Range("A1:G1").Select() | Out-Null

# --- New sheet: ImageLocators (inserted right after TableLocators) ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$imageSheet = $wb.Worksheets.Add([System.Type]::Missing, $afterSheet)
$imageSheet.Name = "ImageLocators"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $imageSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$imageSheet.Range("A1:G1").Select() | Out-Null

# --- New sheet: HeadingLocators (inserted right after ImageLocators) ---
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$headingSheet = $wb.Worksheets.Add([System.Type]::Missing, $afterSheet2)
$headingSheet.Name = "HeadingLocators"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $headingSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
# HeadingLocators ends up the active sheet/tab, with B2 selected.
$headingSheet.Range("B2").Select() | Out-Null
